$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (cell styles) from row 3 into the new row 4 first, so the
# shared-string cells land with the same style indices (s="4" for the number
# column, s="5" for the wrapped-text columns) as the rest of the table.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("C3:E3").Copy()
$ws.Range("C4:E4").PasteSpecial(-4122)

# New data row: row number 291, English source line, Russian translation and
# the "converted" (re-encoded) string, mirroring the existing rows.
$ws.Range("B4").Value = 291
$ws.Range("C4").Value = ' How was the expedition? Was it\nfun enough for you? Ha ha ha!'
$ws.Range("D4").Value = 'Как прошла экспедиция? Вам\nбыло весело? Ха-ха-ха!'
$ws.Range("E4").Value = ' Ëàë ðñïšìà üëòðåäéøéÿ? Âàí\náúìï âåòåìï? Öà-öà-öà!'

$ws.Range("A4:E4").RowHeight = 21.6

[void]$ws.Range("C2").Select()
